# Update data parsing logic
# Appends a new row (row 86) of parsed data to each of the four sheets
# (FE_LFT_#1, FE_LFT_#2, FE_PLT_#1, FE_PLT_#2), mirroring the structure
# of the existing rows (time, total length, ID, actual length, checksum
# and their decimal equivalents).

$wb = $excel.ActiveWorkbook

$newRows = @{
    "FE_LFT_#1" = @{
        Row = 86
        A = 45872.49319444445
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x14"
        E = "0xf"
        F = 380
        G = 759863127514710900000000.0
        H = 276
        I = 15
    }
    "FE_LFT_#2" = @{
        Row = 86
        A = 45872.49319444445
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x24"
        E = "0xe"
        F = 400
        G = 568432987514711000000000.0
        H = 292
        I = 14
    }
    "FE_PLT_#1" = @{
        Row = 86
        A = 45872.49319444445
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x5F"
        E = "0x3"
        F = 110
        G = 568631262647114000000000.0
        H = 95
        I = 3
    }
    "FE_PLT_#2" = @{
        Row = 86
        A = 45872.49319444445
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x5D"
        E = "0x3"
        F = 110
        G = 985046333984776000000000.0
        H = 93
        I = 3
    }
}

foreach ($sheetName in $newRows.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $data = $newRows[$sheetName]
    $r = $data.Row

    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E

    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
}
